$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect to allow editing, then restore protection at the end.
$ws.Unprotect()

# Update the confidential disclosure date string (A59)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-56
$ws.Range("D2").Value = 0.01533871547218321
$ws.Range("E2").Value = -0.002921445574550985
$ws.Range("D3").Value = 0.05058335661211796
$ws.Range("E3").Value = -0.03072052885973731
$ws.Range("D4").Value = 0.01439767980103758
$ws.Range("E4").Value = -0.01063080371437708
$ws.Range("D5").Value = 0.009600297283919521
$ws.Range("E5").Value = 0.02253809706748622
$ws.Range("D6").Value = 0.01521073609617398
$ws.Range("E6").Value = 0.007969222313822355
$ws.Range("D7").Value = 0.01952284811189684
$ws.Range("E7").Value = 0.01864594894561589
$ws.Range("D8").Value = 0.004674843186630845
$ws.Range("E8").Value = -0.02324083022688572
$ws.Range("D9").Value = 0.006907566954688002
$ws.Range("E9").Value = 0.002402691013935643
$ws.Range("D10").Value = 0.01450251593253793
$ws.Range("E10").Value = -0.01016612943218431
$ws.Range("D11").Value = 0.008329109245557017
$ws.Range("E11").Value = 0.01033210332103329
$ws.Range("D12").Value = 0.01561938493945416
$ws.Range("E12").Value = -0.002597402597402598
$ws.Range("D13").Value = 0.002878613918593617
$ws.Range("E13").Value = -0.01831091180866973
$ws.Range("D14").Value = 0.005943480243194892
$ws.Range("E14").Value = -0.02673492605233208
$ws.Range("D15").Value = 0.01437143234976192
$ws.Range("E15").Value = -0.0001240387000743803
$ws.Range("D16").Value = 0.01038403265573393
$ws.Range("E16").Value = 0.0006659563132658253
$ws.Range("D17").Value = 0.02070582750252965
$ws.Range("E17").Value = -0.01790127579987977
$ws.Range("D18").Value = 0.008318290623953933
$ws.Range("E18").Value = 0.004034761018001376
$ws.Range("D19").Value = 0.01657216127381677
$ws.Range("E19").Value = 0.01050445103857567
$ws.Range("D20").Value = 0.01201235814588009
$ws.Range("E20").Value = -0.003377341111452159
$ws.Range("D21").Value = 0.007314371714740372
$ws.Range("E21").Value = 0.01461442786069678
$ws.Range("D22").Value = 0.01464865952834097
$ws.Range("E22").Value = 0.008694630478204513
$ws.Range("D23").Value = 0.01977336681839329
$ws.Range("E23").Value = 0.003273459795198885
$ws.Range("D24").Value = 0.01016624642653091
$ws.Range("E24").Value = -0.005435735577766221
$ws.Range("D25").Value = 0.02027483451747265
$ws.Range("E25").Value = -0.0127699826580483
$ws.Range("D26").Value = 0.01403808357161685
$ws.Range("E26").Value = -0.006298850574712689
$ws.Range("D27").Value = 0.02131741770502905
$ws.Range("E27").Value = -0.0699241919624447
$ws.Range("D28").Value = 0.05602755132138638
$ws.Range("E28").Value = -0.02580446970278794
$ws.Range("D29").Value = 0.02134833683381514
$ws.Range("E29").Value = -0.01652173913043464
$ws.Range("D30").Value = 0.02994791161944929
$ws.Range("E30").Value = -0.04433497536945807
$ws.Range("D31").Value = 0.01567609049870896
$ws.Range("E31").Value = -0.06463547334058761
$ws.Range("D32").Value = 0.0133539901635445
$ws.Range("E32").Value = -0.0003498338289313896
$ws.Range("D33").Value = 0.01902651311113485
$ws.Range("E33").Value = -0.05966503838101889
$ws.Range("D34").Value = 0.04337154665956364
$ws.Range("E34").Value = -0.02558749622650325
$ws.Range("D35").Value = 0.01078788688262182
$ws.Range("E35").Value = 0.00786324786324788
$ws.Range("D36").Value = 0.009961122902158779
$ws.Range("E36").Value = 0.004060475161987043
$ws.Range("D37").Value = 0.01062883470422419
$ws.Range("E37").Value = -0.01474734330947736
$ws.Range("D38").Value = 0.007267532000754288
$ws.Range("E38").Value = 0
$ws.Range("D39").Value = 0.01205643173502448
$ws.Range("E39").Value = -0.003082029397818875
$ws.Range("D40").Value = 0.01740685481215934
$ws.Range("E40").Value = -0.005053340819764163
$ws.Range("D41").Value = 0.01722914665844048
$ws.Range("E41").Value = -0.0006778741865511861
$ws.Range("D42").Value = 0.03194874192161077
$ws.Range("E42").Value = -0.02274170274170273
$ws.Range("D43").Value = 0.0113964171363519
$ws.Range("E43").Value = -0.0009797210049724558
$ws.Range("D44").Value = 0.02211584427322492
$ws.Range("E44").Value = -0.02649491642254009
$ws.Range("D45").Value = 0.01251892780855637
$ws.Range("E45").Value = -0.03641345176543365
$ws.Range("D46").Value = 0.008704718264339129
$ws.Range("E46").Value = -0.01722329911976872
$ws.Range("D47").Value = 0.01333216851201555
$ws.Range("E47").Value = 0.006722270633636507
$ws.Range("D48").Value = 0.01037047864400961
$ws.Range("E48").Value = -0.0007201728414820696
$ws.Range("D49").Value = 0.01611415246947252
$ws.Range("E49").Value = -0.01874506482954341
$ws.Range("D50").Value = 0.008538013140546134
$ws.Range("E50").Value = 0.001529894131326204
$ws.Range("D51").Value = 0.01213117857519125
$ws.Range("E51").Value = 0.005583902955617503
$ws.Range("D52").Value = 0.008093343204882966
$ws.Range("E52").Value = 0.03024353565311833
$ws.Range("D53").Value = 0.01000587265515336
$ws.Range("E53").Value = 0.01051432011696907
$ws.Range("D54").Value = 0.1334651547952782
$ws.Range("E54").Value = -0.0001971220185293943
$ws.Range("D55").Value = 0.04376900806459517
$ws.Range("E55").Value = -0.00907246030798603
$ws.Range("E56").Value = -0.01194785457334935

# Restore sheet protection to its prior state
$ws.Protect()
